# Replace the static "www.drpaulduenas.com" text in the primary footer
# with a MERGEFIELD complex field ("=website") so the site URL becomes
# configurable, matching the new "=emergency_number" / branch-office style
# fields already used elsewhere in the footer.

$d = $word.ActiveDocument

# The default ("primary") footer - the one that actually holds the
# "www.drpaulduenas.com" run - is Footers.Item(1) (wdHeaderFooterPrimary).
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)
$rng = $footer.Range

$found = $rng.Find.Execute("www.drpaulduenas.com", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'www.drpaulduenas.com' in the primary footer"
}

# Collapse the found range to nothing so InsertXML replaces it in place
# rather than appending after it.
$rng.Text = ""

# Run properties shared by every run that made up the old text and by
# every run of the new field (rFonts Avenir Book, bold, sz 20).
$runPr = '<w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'

# Rebuild the owning paragraph verbatim (same pPr / paragraph identity)
# but with the single text run swapped out for a MERGEFIELD complex
# field: begin -> instrText -> separate -> cached result -> end.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
'<w:body>' +
'<w:p w14:paraId="36E3463E" w14:textId="77777777" w:rsidR="00EC5285" w:rsidRDefault="00EC5285" w:rsidP="00EC5285">' +
'<w:pPr><w:pStyle w:val="Footer"/><w:jc w:val="center"/>' + $runPr + '</w:pPr>' +
'<w:r>' + $runPr + '<w:fldChar w:fldCharType="begin"/></w:r>' +
'<w:r>' + $runPr + '<w:instrText xml:space="preserve"> MERGEFIELD =website \* MERGEFORMAT </w:instrText></w:r>' +
'<w:r>' + $runPr + '<w:fldChar w:fldCharType="separate"/></w:r>' +
'<w:r>' + $runPr + '<w:t>&#171;=website&#187;</w:t></w:r>' +
'<w:r>' + $runPr + '<w:fldChar w:fldCharType="end"/></w:r>' +
'</w:p>' +
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml)
